# Edit the "Liste-referentiels" workbook:
#  - rename the "RDF: classes utilisees" column header
#  - update the "agents/producteurs" row with the new EAC-CPF/RDF file description
#  - enable word-wrap on column A (file/folder name column)
#  - grow row 6 to fit the new, longer description
#  - move the selection to G17 (where the new description now lives)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: rename "RDF: classes utilisees" -> "RDF: principales classes utilisees"
$ws.Range("D3").Value = "RDF: principales classes utilisées"

# --- Row 6 ("agents/producteurs") content updates
$ws.Range("A6").Value = "agents/producteurs/eac-cpf et agents/producteurs/rdf"

$ws.Range("D6").Value = "rico:Person; rico:Family; rico:CorporateBody; rico:CorporateBodyType; rico:OccupationType; rico:ActivityType; rico:Relation et ses sous-classes ; rico:Place"

$ws.Range("G6").Value = "Référentiel produit à partir du référentiel des producteurs du SIA, en utilisant le logiciel RiC-O Converter (https://github.com/ArchivesNationalesFR/rico-converter), et en rendant ensuite par script le résultat conforme à RiC-O 0.2. Un fichier par agent, 6 fichiers pour les relations entre agents et relations de provenance entre agents et documents, + un fichier pour les lieux (qui sera repris prochainement). Fourni avec une liste des notices EAC-CPF au format tsv (en utf-8; séparateur : tabulation). ATTENTION : les notices dont la liste TSV indique qu'elles ne contiennent pas d'éléments biographiques ou historiques rédigés (pas d'élément EAC-CPF biogHist) et/ou qu'elles ont été créées en 2013 ne sont pas forcément fiables (des vérifications et enrichissements restent à faire) ; leur version RDF ne l'est pas plus."

# Row 6 grew considerably taller because of the new, much longer description
$ws.Rows(6).RowHeight = 180

# --- Enable word-wrap for the "Nom du fichier ou du dossier" column (A), like the other text columns
$ws.Range("A3:A19").WrapText = $true

# --- Update the saved selection/view to reflect where the editor ended up
$ws.Range("G17").Select()
